$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C52").Value = "'249"
$ws.Range("D52").Value = "'876892.59"
$ws.Range("C52").Style = "Normal"
$ws.Range("D52").Style = "Normal"
$ws.Range("C56").Value = "'35"
$ws.Range("D56").Value = "'82500.00"
$ws.Range("C56").Style = "Normal"
$ws.Range("D56").Style = "Normal"
$ws.Range("C58").Value = "'22"
$ws.Range("D58").Value = "'53000.00"
$ws.Range("C58").Style = "Normal"
$ws.Range("D58").Style = "Normal"
$ws.Range("C60").Value = "'37"
$ws.Range("D60").Value = "'146456.00"
$ws.Range("C60").Style = "Normal"
$ws.Range("D60").Style = "Normal"
$ws.Range("C61").Value = "'64"
$ws.Range("D61").Value = "'138500.00"
$ws.Range("C61").Style = "Normal"
$ws.Range("D61").Style = "Normal"
$ws.Range("C74").Value = "'12"
$ws.Range("D74").Value = "'35000.00"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").Style = "Normal"
$ws.Range("C77").Value = "'80"
$ws.Range("D77").Value = "'210987.00"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").Style = "Normal"
$ws.Range("C78").Value = "'181"
$ws.Range("D78").Value = "'487193.00"
$ws.Range("C78").Style = "Normal"
$ws.Range("D78").Style = "Normal"
$ws.Range("C80").Value = "'420"
$ws.Range("D80").Value = "'1665769.58"
$ws.Range("C80").Style = "Normal"
$ws.Range("D80").Style = "Normal"
$ws.Range("C89").Value = "'96"
$ws.Range("D89").Value = "'231510.00"
$ws.Range("C89").Style = "Normal"
$ws.Range("D89").Style = "Normal"
$ws.Range("C122").Value = "'220"
$ws.Range("D122").Value = "'599708.00"
$ws.Range("C122").Style = "Normal"
$ws.Range("D122").Style = "Normal"
$ws.Range("C123").Value = "'70"
$ws.Range("D123").Value = "'196976.45"
$ws.Range("C123").Style = "Normal"
$ws.Range("D123").Style = "Normal"
$ws.Range("C124").Value = "'414"
$ws.Range("D124").Value = "'1712502.18"
$ws.Range("C124").Style = "Normal"
$ws.Range("D124").Style = "Normal"
$ws.Range("C132").Value = "'68"
$ws.Range("D132").Value = "'272186.75"
$ws.Range("C132").Style = "Normal"
$ws.Range("D132").Style = "Normal"
$ws.Range("C133").Value = "'110"
$ws.Range("D133").Value = "'276626.44"
$ws.Range("C133").Style = "Normal"
$ws.Range("D133").Style = "Normal"
$ws.Range("C191").Value = "'109"
$ws.Range("D191").Value = "'290000.00"
$ws.Range("C191").Style = "Normal"
$ws.Range("D191").Style = "Normal"
$ws.Range("C194").Value = "'594"
$ws.Range("D194").Value = "'2118871.11"
$ws.Range("C194").Style = "Normal"
$ws.Range("D194").Style = "Normal"
$ws.Range("C203").Value = "'127"
$ws.Range("D203").Value = "'281196.77"
$ws.Range("C203").Style = "Normal"
$ws.Range("D203").Style = "Normal"
$ws.Range("C222").Value = "'157"
$ws.Range("D222").Value = "'429905.00"
$ws.Range("C222").Style = "Normal"
$ws.Range("D222").Style = "Normal"
$ws.Range("C223").Value = "'16"
$ws.Range("D223").Value = "'48000.00"
$ws.Range("C223").Style = "Normal"
$ws.Range("D223").Style = "Normal"
$ws.Range("C224").Value = "'314"
$ws.Range("D224").Value = "'1049983.50"
$ws.Range("C224").Style = "Normal"
$ws.Range("D224").Style = "Normal"
$ws.Range("C228").Value = "'78"
$ws.Range("D228").Value = "'227187.09"
$ws.Range("C228").Style = "Normal"
$ws.Range("D228").Style = "Normal"
$ws.Range("C229").Value = "'27"
$ws.Range("D229").Value = "'82587.00"
$ws.Range("C229").Style = "Normal"
$ws.Range("D229").Style = "Normal"
$ws.Range("C230").Value = "'26"
$ws.Range("D230").Value = "'74600.00"
$ws.Range("C230").Style = "Normal"
$ws.Range("D230").Style = "Normal"
$ws.Range("C232").Value = "'52"
$ws.Range("D232").Value = "'177270.00"
$ws.Range("C232").Style = "Normal"
$ws.Range("D232").Style = "Normal"
